$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.475.51"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "3.151.89"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'593.22"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").Value = "'146.98"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.148.43"
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("D11").Value = "'5.94"
$ws.Range("E11").Value = "  +3.64%  "
$ws.Range("D12").Value = "'0.464"
$ws.Range("E12").Value = "  -0.86%  "
$ws.Range("E13").Value = "  -1.92%  "
$ws.Range("D14").Value = "'37.47"
$ws.Range("E14").Value = "  +3.37%  "
$ws.Range("D15").Value = "3.673.40"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").Value = "'7.28"
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("D18").Value = "64.248.64"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").Value = "3.155.10"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("D20").Value = "'469.46"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").Value = "'14.49"
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("D22").Value = "'0.737"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "'13.16"
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").Value = "'2.33"
$ws.Range("E25").Value = "  +5.87%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").Value = "'81.72"
$ws.Range("E26").Value = "  -0.91%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "'9.66"
$ws.Range("E28").Value = "  +7.90%  "
$ws.Range("D29").Value = "'7.50"
$ws.Range("E29").Value = "  +8.68%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'2.27"
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.74"
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").Value = "'27.54"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("D35").Value = "0.0₃0848"
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("D37").Value = "'6.27"
$ws.Range("E37").Value = "  +2.73%  "
$ws.Range("D38").Value = "'2.34"
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("E39").Value = "  -3.79%  "
$ws.Range("D40").Value = "'51.92"
$ws.Range("E40").Value = "  +1.78%  "
$ws.Range("D41").Value = "'458.60"
$ws.Range("E41").Value = "  +1.84%  "
$ws.Range("D42").Value = "'9.28"
$ws.Range("E42").Value = "  +5.59%  "
$ws.Range("E43").Value = "  +7.49%  "
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").Value = "2.944.68"
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("D46").Value = "'39.91"
$ws.Range("E46").Value = "  +11.91%  "
$ws.Range("D47").Value = "'0.109"
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("D48").Value = "'129.54"
$ws.Range("E48").Value = "  +2.59%  "
$ws.Range("D50").Value = "'2.25"
$ws.Range("E50").Value = "  +2.84%  "
$ws.Range("E51").Value = "  -0.50%  "
